$d = $word.ActiveDocument

# The "Discipline responsable" row in the header table currently prints the
# Jinja-style placeholder "{{ programme.discipline }}". Switch it over to
# "{{ programme.departement }}" so the generated "plan de cours" pulls the
# department instead of the discipline.
$d.Content.Find.Execute("discipline", $true, $false, $false, $false, $false, $true, 1, $false, "departement", 2)
